# Updates monthly "Inscripcion Registro Automotor" stats (rows 2-13) from the
# 2014 sample data to the corresponding 2022 (Jan-Dec) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 44562
$ws.Range("B2").Value = 70
$ws.Range("C2").Value = 247
$ws.Range("D2").Value = 46
$ws.Range("E2").Value = 124
$ws.Range("F2").Value = 49
$ws.Range("G2").Value = 38
$ws.Range("H2").Value = 46
$ws.Range("I2").Value = 46
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 243
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 47
$ws.Range("O2").Value = 1020

# Row 3
$ws.Range("A3").Value = 44593
$ws.Range("B3").Value = 91
$ws.Range("C3").Value = 317
$ws.Range("D3").Value = 54
$ws.Range("E3").Value = 132
$ws.Range("F3").Value = 58
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = 64
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 36
$ws.Range("K3").Value = 34
$ws.Range("L3").Value = 259
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 72
$ws.Range("O3").Value = 1265

# Row 4
$ws.Range("A4").Value = 44621
$ws.Range("B4").Value = 118
$ws.Range("C4").Value = 378
$ws.Range("D4").Value = 79
$ws.Range("E4").Value = 229
$ws.Range("F4").Value = 87
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = 113
$ws.Range("I4").Value = 58
$ws.Range("J4").Value = 62
$ws.Range("K4").Value = 52
$ws.Range("L4").Value = 363
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = 79
$ws.Range("O4").Value = 1702

# Row 5
$ws.Range("A5").Value = 44652
$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 314
$ws.Range("D5").Value = 115
$ws.Range("E5").Value = 146
$ws.Range("F5").Value = 55
$ws.Range("G5").Value = 57
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 41
$ws.Range("K5").Value = 43
$ws.Range("L5").Value = 269
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = 57
$ws.Range("O5").Value = 1325

# Row 6
$ws.Range("A6").Value = 44682
$ws.Range("B6").Value = 91
$ws.Range("C6").Value = 303
$ws.Range("D6").Value = 87
$ws.Range("E6").Value = 195
$ws.Range("F6").Value = 51
$ws.Range("G6").Value = 58
$ws.Range("H6").Value = 80
$ws.Range("I6").Value = 68
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 26
$ws.Range("L6").Value = 270
$ws.Range("M6").Value = 12
$ws.Range("N6").Value = 69
$ws.Range("O6").Value = 1345

# Row 7
$ws.Range("A7").Value = 44713
$ws.Range("B7").Value = 97
$ws.Range("C7").Value = 345
$ws.Range("D7").Value = 47
$ws.Range("E7").Value = 159
$ws.Range("F7").Value = 47
$ws.Range("G7").Value = 39
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 53
$ws.Range("J7").Value = 43
$ws.Range("K7").Value = 33
$ws.Range("L7").Value = 303
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 65
$ws.Range("O7").Value = 1306

# Row 8
$ws.Range("A8").Value = 44743
$ws.Range("B8").Value = 110
$ws.Range("C8").Value = 288
$ws.Range("D8").Value = 42
$ws.Range("E8").Value = 168
$ws.Range("F8").Value = 36
$ws.Range("G8").Value = 42
$ws.Range("H8").Value = 69
$ws.Range("I8").Value = 58
$ws.Range("J8").Value = 44
$ws.Range("K8").Value = 32
$ws.Range("L8").Value = 263
$ws.Range("M8").Value = 19
$ws.Range("N8").Value = 51
$ws.Range("O8").Value = 1222

# Row 9
$ws.Range("A9").Value = 44774
$ws.Range("B9").Value = 107
$ws.Range("C9").Value = 336
$ws.Range("D9").Value = 36
$ws.Range("E9").Value = 148
$ws.Range("F9").Value = 59
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = 65
$ws.Range("I9").Value = 37
$ws.Range("J9").Value = 38
$ws.Range("K9").Value = 28
$ws.Range("L9").Value = 292
$ws.Range("M9").Value = 9
$ws.Range("N9").Value = 75
$ws.Range("O9").Value = 1270

# Row 10
$ws.Range("A10").Value = 44805
$ws.Range("B10").Value = 92
$ws.Range("C10").Value = 292
$ws.Range("D10").Value = 33
$ws.Range("E10").Value = 140
$ws.Range("F10").Value = 49
$ws.Range("G10").Value = 36
$ws.Range("H10").Value = 65
$ws.Range("I10").Value = 48
$ws.Range("J10").Value = 55
$ws.Range("K10").Value = 37
$ws.Range("L10").Value = 206
$ws.Range("M10").Value = 17
$ws.Range("N10").Value = 66
$ws.Range("O10").Value = 1136

# Row 11
$ws.Range("A11").Value = 44835
$ws.Range("B11").Value = 99
$ws.Range("C11").Value = 258
$ws.Range("D11").Value = 29
$ws.Range("E11").Value = 141
$ws.Range("F11").Value = 44
$ws.Range("G11").Value = 61
$ws.Range("H11").Value = 53
$ws.Range("I11").Value = 40
$ws.Range("J11").Value = 49
$ws.Range("K11").Value = 33
$ws.Range("L11").Value = 266
$ws.Range("M11").Value = 18
$ws.Range("N11").Value = 79
$ws.Range("O11").Value = 1170

# Row 12
$ws.Range("A12").Value = 44866
$ws.Range("B12").Value = 115
$ws.Range("C12").Value = 265
$ws.Range("D12").Value = 37
$ws.Range("E12").Value = 135
$ws.Range("F12").Value = 48
$ws.Range("G12").Value = 42
$ws.Range("H12").Value = 49
$ws.Range("I12").Value = 32
$ws.Range("J12").Value = 36
$ws.Range("K12").Value = 22
$ws.Range("L12").Value = 257
$ws.Range("M12").Value = 9
$ws.Range("N12").Value = 58
$ws.Range("O12").Value = 1105

# Row 13
$ws.Range("A13").Value = 44896
$ws.Range("B13").Value = 63
$ws.Range("C13").Value = 197
$ws.Range("D13").Value = 27
$ws.Range("E13").Value = 111
$ws.Range("F13").Value = 26
$ws.Range("G13").Value = 28
$ws.Range("H13").Value = 10
$ws.Range("I13").Value = 33
$ws.Range("J13").Value = 27
$ws.Range("K13").Value = 19
$ws.Range("L13").Value = 154
$ws.Range("M13").Value = 10
$ws.Range("N13").Value = 40
$ws.Range("O13").Value = 745

Write-Output "Updated rows 2-13 with new monthly registration data"
